# Refresh the cryptocurrency price / 1h-volume table on the active sheet
# (columns B:Coin, C:Link, D:Price, E:Volume(1h), rows 2-51) with newly
# scraped values. A few coins also changed rank/row position.
#
# Cells are briefly switched to Text number format while the value is
# assigned so numeric-looking strings (e.g. "537.37", "1.975.95") land
# as literal text instead of being auto-converted to numbers by Excel,
# then the style is reset to Normal so no formatting residue is left
# behind -- matching how the sheet already stores these as plain text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
  @{Row=2; D="60.024.26"; E="  +1.51%  "},
  @{Row=3; D="2.677.02"; E="  +3.15%  "},
  @{Row=4; E="  +0.20%  "},
  @{Row=5; D="537.37"; E="  +1.48%  "},
  @{Row=6; D="146.17"; E="  +4.26%  "},
  @{Row=7; D="0.999"; E="  +0.15%  "},
  @{Row=8; D="0.575"; E="  +1.49%  "},
  @{Row=9; D="2.678.05"; E="  +2.61%  "},
  @{Row=10; D="6.65"; E="  +3.08%  "},
  @{Row=11; E="  +2.43%  "},
  @{Row=12; D="0.339"; E="  +1.65%  "},
  @{Row=13; E="  -1.21%  "},
  @{Row=14; D="3.137.18"; E="  +2.66%  "},
  @{Row=15; D="59.953.36"; E="  +1.45%  "},
  @{Row=16; D="21.21"; E="  +3.42%  "},
  @{Row=17; D="2.695.65"; E="  +4.42%  "},
  @{Row=18; D="0.0000136"; E="  +1.47%  "},
  @{Row=19; D="345.47"; E="  -0.46%  "},
  @{Row=20; E="  +2.27%  "},
  @{Row=21; D="10.31"; E="  +1.94%  "},
  @{Row=22; D="6.41"; E="  -0.59%  "},
  @{Row=23; D="0.999"; E="  +0.01%  "},
  @{Row=24; D="67.66"; E="  +0.47%  "},
  @{Row=25; D="0.417"; E="  +2.69%  "},
  @{Row=26; D="0.167"; E="  -0.22%  "},
  @{Row=27; D="0.998"; E="  +0.00%  "},
  @{Row=28; D="7.33"; E="  +2.37%  "},
  @{Row=29; D="0.0₃0757"; E="  +2.48%  "},
  @{Row=30; E="  +0.04%  "},
  @{Row=31; D="1.67"; E="  +3.03%  "},
  @{Row=32; E="  +0.82%  "},
  @{Row=33; D="19.17"; E="  +1.93%  "},
  @{Row=34; D="150.55"; E="  +0.89%  "},
  @{Row=35; D="4.06"; E="  +1.68%  "},
  @{Row=36; D="1.16"; E="  +2.80%  "},
  @{Row=37; D="1.48"; E="  +0.19%  "},
  @{Row=38; D="0.843"; E="  +1.55%  "},
  @{Row=39; D="0.829"; E="  +0.86%  "},
  @{Row=40; D="293.69"; E="  +9.10%  "},
  @{Row=41; D="3.62"; E="  +2.23%  "},
  @{Row=42; E="  +0.23%  "},
  @{Row=43; D="0.607"; E="  +1.50%  "},
  @{Row=44; B="WhiteBITCoin"; C="https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"; D="10.75"; E="  +0.01%  "},
  @{Row=45; B="Stellar"; C="https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"; D="0.0958"; E="  -0.04%  "},
  @{Row=46; B="Hedera"; C="https://coinranking.com/coin/jad286TjB+hedera-hbar"; D="0.0539"; E="  +3.80%  "},
  @{Row=47; B="VeChain"; C="https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; D="0.0227"; E="  +2.35%  "},
  @{Row=48; B="Maker"; C="https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"; D="1.975.95"; E="  +0.64%  "},
  @{Row=49; D="4.62"; E="  -0.04%  "},
  @{Row=50; D="18.53"; E="  +1.32%  "},
  @{Row=51; D="110.07"; E="  -0.72%  "}
)

function Set-TextCell($cell, $value) {
  $cell.NumberFormat = "@"
  $cell.Value = $value
  $cell.Style = "Normal"
}

foreach ($u in $updates) {
  if ($u.ContainsKey("B")) { Set-TextCell $ws.Cells.Item($u.Row, 2) $u.B }
  if ($u.ContainsKey("C")) { Set-TextCell $ws.Cells.Item($u.Row, 3) $u.C }
  if ($u.ContainsKey("D")) { Set-TextCell $ws.Cells.Item($u.Row, 4) $u.D }
  if ($u.ContainsKey("E")) { Set-TextCell $ws.Cells.Item($u.Row, 5) $u.E }
}
